$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Date of Birth and Guardian Phone Number become plain text values ---
# Switch number format to Text ("@") first so the literal strings below are
# not re-interpreted as a date serial / number by Excel's auto-detection.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "01-01-2000"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "1234567890"

# --- Row 3: only keep a blank, text-formatted E3; drop the rest of the row ---
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").ClearContents()

$ws.Range("A3").Clear()
$ws.Range("B3").Clear()
$ws.Range("F3").Clear()
$ws.Range("G3").Clear()

# --- Row 4: only keep a blank, text-formatted I4; drop the rest of the row ---
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").ClearContents()

$ws.Range("A4").Clear()
$ws.Range("B4").Clear()
$ws.Range("D4").Clear()
$ws.Range("F4").Clear()
$ws.Range("G4").Clear()
